$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing row 2 values ---
$ws.Range("A2").Value = "wiGka687"
$ws.Range("B2").Value = 23110861
$ws.Range("C2").Value = "mciswee13"
$ws.Range("D2").Value = "S`$e5w%V9"
$ws.Range("E2").Value = "MR"
$ws.Range("F2").Value = "zpFXWIVn"
$ws.Range("G2").Value = "RVgJ"
$ws.Range("H2").Value = "Candidate"

# --- Add new row 3 with the same visual style as row 2 (bordered data row) ---
$r3 = $ws.Range("A3:H3")
$r3.Style = "Normal"
$r3.Borders.Color = 0
$r3.Borders.Weight = 2
$r3.Borders.LineStyle = 1

$ws.Range("A3").Value = "iRiaL673"
$ws.Range("B3").Value = 23110860
$ws.Range("C3").Value = "zwqhjhs73"
$ws.Range("D3").Value = "X#4meH2`$"
$ws.Range("E3").Value = "MR"
$ws.Range("F3").Value = "psabtdSX"
$ws.Range("G3").Value = "TVAM"
$ws.Range("H3").Value = "Candidate"

# --- Update the sheet view selection to span the new range ---
$null = $ws.Range("A1:H3").Select()
